$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.779.47"
$ws.Range("E2").Value = "  +0.00%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.279.31"
$ws.Range("E3").Value = "  +1.21%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.14"
$ws.Range("E5").Value = "  +0.81%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  +1.61%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.29"
$ws.Range("E7").Value = "  +6.82%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.646"
$ws.Range("E9").Value = "  -2.94%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.47"
$ws.Range("E10").Value = "  +1.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0983"
$ws.Range("E11").Value = "  +2.69%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.46"
$ws.Range("E12").Value = "  -0.30%  "

# Row 13
$ws.Range("E13").Value = "  +2.07%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.619.64"
$ws.Range("E14").Value = "  +1.42%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.02"
$ws.Range("E15").Value = "  +1.00%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.872"
$ws.Range("E16").Value = "  -0.77%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.277.81"
$ws.Range("E17").Value = "  +1.23%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.658.68"
$ws.Range("E18").Value = "  -0.12%  "

# Row 19
$ws.Range("E19").Value = "  +1.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.22"
$ws.Range("E20").Value = "  -1.16%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.60"
$ws.Range("E21").Value = "  -0.51%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.14"
$ws.Range("E22").Value = "  -0.12%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.16"
$ws.Range("E23").Value = "  +5.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.86"
$ws.Range("E24").Value = "  -1.80%  "

# Row 25
$ws.Range("E25").Value = "  -0.05%  "

# Row 26
$ws.Range("E26").Value = "  -1.56%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.39"
$ws.Range("E27").Value = "  -1.22%  "

# Row 28
$ws.Range("E28").Value = "  +2.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.80"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.08"
$ws.Range("E30").Value = "  -1.60%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0880"
$ws.Range("E31").Value = "  +9.92%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.46"
$ws.Range("E32").Value = "  -0.44%  "

# Row 33
$ws.Range("E33").Value = "  +0.88%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.46"
$ws.Range("E34").Value = "  +0.44%  "

# Row 35
$ws.Range("E35").Value = "  +1.97%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.52"
$ws.Range("E36").Value = "  +2.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.78"
$ws.Range("E37").Value = "  +1.68%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0305"
$ws.Range("E38").Value = "  -5.43%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.74"
$ws.Range("E39").Value = "  +10.68%  "

# Row 40
$ws.Range("E40").Value = "  -0.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.89"
$ws.Range("E41").Value = "  +1.47%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.209"
$ws.Range("E42").Value = "  +3.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "61.42"
$ws.Range("E43").Value = "  -1.37%  "

# Row 44
$ws.Range("E44").Value = "  -0.22%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.54"
$ws.Range("E45").Value = "  +11.16%  "

# Row 46
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.75"
$ws.Range("E46").Value = "  -2.85%  "

# Row 47
$ws.Range("E47").Value = "  -0.90%  "

# Row 48
$ws.Range("E48").Value = "  +0.09%  "

# Row 49
$ws.Range("E49").Value = "  +0.73%  "

# Row 50
$ws.Range("E50").Value = "  -1.37%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.24"
$ws.Range("E51").Value = "  -0.98%  "
